$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("May")

# --- Row 18: problem 378 (kth-smallest-element-in-a-sorted-matrix) ---
$ws.Range("B18").Value = 378
$ws.Range("D18").Value = "Medium"
$ws.Range("E18").Value = "Yes"
$ws.Range("F18").Value = "Yes"
$ws.Range("G18").Value = "No"
$ws.Range("H18").Value = "No"
$ws.Range("I18").Value = "No"
$ws.Range("J18").Value = "No"
$ws.Range("K18").Value = "No"

# --- Row 24: problem 22 (generate-parentheses) ---
$ws.Range("B24").Value = 22
$ws.Range("C24").Value = "https://leetcode.com/problems/generate-parentheses/"
$ws.Range("D24").Value = "Medium"
$ws.Range("E24").Value = "Yes"
$ws.Range("F24").Value = "Yes"
$ws.Range("G24").Value = "Yes"
$ws.Range("H24").Value = "No"
$ws.Range("I24").Value = "No"
$ws.Range("J24").Value = "No"
$ws.Range("K24").Value = "No"

# --- Row 25: problem 33 (search-in-rotated-sorted-array) ---
$ws.Range("B25").Value = 33
$ws.Range("C25").Value = "https://leetcode.com/problems/search-in-rotated-sorted-array/"
$ws.Range("D25").Value = "Medium"
$ws.Range("E25").Value = "Yes"
$ws.Range("F25").Value = "Yes"
$ws.Range("G25").Value = "Yes"
$ws.Range("H25").Value = "No"
$ws.Range("I25").Value = "Yes"
$ws.Range("J25").Value = "No"
$ws.Range("K25").Value = "Yes"

# --- Convert A1:K62 into a proper Excel Table with the Light20 style ---
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:K62"), $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = "TableStyleLight20"

# --- Move the active selection from K19 to B19 ---
$ws.Range("B19").Select()
